$d = $word.ActiveDocument

# The document starts with a content control (Structured Document Tag)
# wrapping the auto-generated "Table of Contents" building block: a
# "TOC Heading" paragraph ("Table of Contents") followed by a paragraph
# holding the TOC field codes. The edited document no longer has a Table
# of Contents, so remove the whole block -- the sdt wrapper together with
# both of its paragraphs.
if ($d.ContentControls.Count -gt 0) {
    $cc = $d.ContentControls.Item(1)

    # Make sure it is not locked so it (and its contents) can be removed.
    $cc.LockContentControl = $false
    $cc.LockContents = $false

    $start = $cc.Range.Start
    $paraIndex = 1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Start -eq $start) {
            $paraIndex = $i
            break
        }
    }

    # Removing the content control itself only strips the sdt wrapper and
    # leaves its paragraphs (the TOC heading text + the field-code
    # paragraph) behind in the body, as two separate paragraphs right
    # where the control used to be.
    $cc.Delete()

    # The former TOC block is exactly those two paragraphs (heading text
    # + field-code paragraph); grab their combined range and delete it
    # (including both paragraph marks) so nothing is left behind.
    $headingPara = $d.Paragraphs($paraIndex)
    $fieldPara = $d.Paragraphs($paraIndex + 1)
    $rng = $d.Range($headingPara.Range.Start, $fieldPara.Range.End)
    $rng.Delete()
}
